$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Year -> Data value map reflecting the revised GDP per Capita series.
# Existing rows (years 1950-2008, sheet rows 2-60) get their Data column
# (E) revised in place; new rows are appended for years 2009-2016.
$values = @{
    1950 = "815"
    1951 = "834"
    1952 = "851"
    1953 = "870"
    1954 = "888"
    1955 = "907"
    1956 = "926"
    1957 = "945"
    1958 = "966"
    1959 = "983"
    1960 = "974"
    1961 = "971"
    1962 = "1004"
    1963 = "972"
    1964 = "929"
    1965 = "913"
    1966 = "878"
    1967 = "866"
    1968 = "843"
    1969 = "881"
    1970 = "880"
    1971 = "877"
    1972 = "800"
    1973 = "746"
    1974 = "829"
    1975 = "952"
    1976 = "918"
    1977 = "830"
    1978 = "808"
    1979 = "622"
    1980 = "583"
    1981 = "521"
    1982 = "524"
    1983 = "577"
    1984 = "601"
    1985 = "642"
    1986 = "662"
    1987 = "663"
    1988 = "692"
    1989 = "685"
    1990 = "692"
    1991 = "703.283742270115"
    1992 = "698.739002074788"
    1993 = "680.800540476404"
    1994 = "713.827084771773"
    1995 = "700.707623511954"
    1996 = "694.105019824845"
    1997 = "719.153116123223"
    1998 = "753.860799416204"
    1999 = "733.781232654736"
    2000 = "716.042363333227"
    2001 = "783.055164226115"
    2002 = "831.829185503164"
    2003 = "925.512616305639"
    2004 = "1200.27255858474"
    2005 = "1266.10503289822"
    2006 = "1225.51631209176"
    2007 = "1422.34078974323"
    2008 = "1422.23120099355"
    2009 = "1533.17948604512"
    2010 = "1720.33434710546"
    2011 = "1661"
    2012 = "1793"
    2013 = "1888"
    2014 = "2279"
    2015 = "2384"
    2016 = "2189"
}

# Revise the Data value for each existing year row (sheet rows 2-60, years 1950-2008)
$row = 2
foreach ($year in 1950..2008) {
    $ws.Cells.Item($row, 5).Value = "'" + $values[[string]$year]
    $row = $row + 1
}

# Append new rows for years 2009-2016 (sheet rows 61-68)
foreach ($year in 2009..2016) {
    $ws.Cells.Item($row, 1).Value = 148
    $ws.Cells.Item($row, 2).Value = "Chad"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year
    $ws.Cells.Item($row, 5).Value = "'" + $values[[string]$year]
    $row = $row + 1
}
